$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 595.7143
$ws.Range("I2").Value = 624.1667
$ws.Range("J2").Value = 425
$ws.Range("K2").Value = 624.1667
$ws.Range("L2").Value = 425
$ws.Range("M2").Value = -511.1667
$ws.Range("N2").Value = -651
# Row 33
$ws.Range("H33").Value = 519.375
$ws.Range("I33").Value = 539.2
$ws.Range("J33").Value = 222
$ws.Range("K33").Value = 539.2
$ws.Range("L33").Value = 222
$ws.Range("M33").Value = -310.2
$ws.Range("N33").Value = -680
# Row 69
$ws.Range("H69").Value = 13999.333
$ws.Range("J69").Value = 15000
$ws.Range("L69").Value = 45000
$ws.Range("N69").Value = -46748
# Row 70
$ws.Range("H70").Value = 6871.4
$ws.Range("I70").Value = 5440.2
$ws.Range("J70").Value = 7587
$ws.Range("K70").Value = 16320.6
$ws.Range("L70").Value = 22761
$ws.Range("M70").Value = -16050.6
$ws.Range("N70").Value = -23301
# Row 72
$ws.Range("H72").Value = 13999.333
$ws.Range("J72").Value = 15000
$ws.Range("L72").Value = 135000
$ws.Range("N72").Value = -143736
# Row 73
$ws.Range("H73").Value = 6871.4
$ws.Range("I73").Value = 5440.2
$ws.Range("J73").Value = 7587
$ws.Range("K73").Value = 16320.6
$ws.Range("L73").Value = 22761
$ws.Range("M73").Value = -15384.6
$ws.Range("N73").Value = -24633
# Row 116
$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
# Row 132
$ws.Range("H132").Value = 5013.5
$ws.Range("I132").Value = 5013.5
$ws.Range("K132").Value = 15040.5
$ws.Range("M132").Value = -12510.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5610.905
$ws.Range("I32").Value = 4364.737
$ws.Range("K32").Value = 4364.737
$ws.Range("M32").Value = -4077.737
# Row 36
$ws.Range("H36").Value = 792.5
$ws.Range("I36").Value = 792.5
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 792.5
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -446.5
$ws.Range("N36").ClearContents()
# Row 43
$ws.Range("H43").Value = 28877.4
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 28877.4
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 28877.4
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -29503.4
# Row 61
$ws.Range("H61").Value = 3115.875
$ws.Range("I61").Value = 2818.6667
$ws.Range("J61").Value = 4007.5
$ws.Range("K61").Value = 2818.6667
$ws.Range("L61").Value = 4007.5
$ws.Range("M61").Value = -2606.6667
$ws.Range("N61").Value = -4431.5
# Row 74
$ws.Range("H74").Value = 1695.8
$ws.Range("I74").Value = 1638.5
$ws.Range("K74").Value = 1638.5
$ws.Range("M74").Value = -764.5
# Row 77
$ws.Range("H77").Value = 1695.8
$ws.Range("I77").Value = 1638.5
$ws.Range("K77").Value = 8192.5
$ws.Range("M77").Value = -3824.5
# Row 122
$ws.Range("H122").Value = 2205.087
$ws.Range("I122").Value = 2168.9546
$ws.Range("K122").Value = 6506.8638
$ws.Range("M122").Value = -4056.8638
# Row 132
$ws.Range("H132").Value = 6514.4443
$ws.Range("I132").Value = 6578.75
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 19736.25
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -17206.25
$ws.Range("N132").Value = -23060
# Row 136
$ws.Range("H136").Value = 3115.875
$ws.Range("I136").Value = 2818.6667
$ws.Range("J136").Value = 4007.5
$ws.Range("K136").Value = 8456.000100000001
$ws.Range("L136").Value = 12022.5
$ws.Range("M136").Value = -5906.000100000001
$ws.Range("N136").Value = -17122.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 352.2857
$ws.Range("I94").Value = 352.2857
$ws.Range("K94").Value = 352.2857
$ws.Range("M94").Value = 98.71429999999998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 1650
$ws.Range("I17").Value = 1600
$ws.Range("K17").Value = 1600
$ws.Range("M17").Value = -1426
# Row 50
$ws.Range("H50").Value = 16085.714
$ws.Range("J50").Value = 31000
$ws.Range("L50").Value = 31000
$ws.Range("N50").Value = -32250
# Row 132
$ws.Range("H132").Value = 3241.8572
$ws.Range("I132").Value = 3313.1667
$ws.Range("J132").Value = 2814
$ws.Range("K132").Value = 9939.500100000001
$ws.Range("L132").Value = 8442
$ws.Range("M132").Value = -7409.500100000001
$ws.Range("N132").Value = -13502

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 115.14286
$ws.Range("J2").Value = 114.166664
$ws.Range("L2").Value = 684.999984
$ws.Range("N2").Value = -910.999984
# Row 12
$ws.Range("H12").Value = 376.8
$ws.Range("I12").Value = 371.33334
$ws.Range("J12").Value = 385
$ws.Range("K12").Value = 1114.00002
$ws.Range("L12").Value = 1155
$ws.Range("M12").Value = -941.0000199999999
$ws.Range("N12").Value = -1501
# Row 26
$ws.Range("H26").Value = 1754.84
$ws.Range("I26").Value = 2000
$ws.Range("J26").Value = 774.2
$ws.Range("K26").Value = 6000
$ws.Range("L26").Value = 2322.6
$ws.Range("M26").Value = -5712
$ws.Range("N26").Value = -2898.6
# Row 38
$ws.Range("H38").Value = 266.4
$ws.Range("J38").Value = 121
$ws.Range("L38").Value = 363
$ws.Range("N38").Value = -1057
# Row 81
$ws.Range("H81").Value = 1694.5
$ws.Range("I81").Value = 1694.5
$ws.Range("K81").Value = 5083.5
$ws.Range("M81").Value = -3960.5
# Row 84
$ws.Range("H84").Value = 1694.5
$ws.Range("I84").Value = 1694.5
$ws.Range("K84").Value = 15250.5
$ws.Range("M84").Value = -9634.5
# Row 97
$ws.Range("H97").Value = 1833.3572
$ws.Range("I97").Value = 1079.75
$ws.Range("J97").Value = 2134.8
$ws.Range("K97").Value = 3239.25
$ws.Range("L97").Value = 6404.400000000001
$ws.Range("M97").Value = -2743.25
$ws.Range("N97").Value = -7396.400000000001
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
# Row 131
$ws.Range("H131").Value = 2644.3333
$ws.Range("J131").Value = 3933
$ws.Range("L131").Value = 11799
$ws.Range("N131").Value = -21879
# Row 132
$ws.Range("H132").Value = 1316.6923
$ws.Range("I132").Value = 1260.6364
$ws.Range("K132").Value = 11345.7276
$ws.Range("M132").Value = -8815.7276

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
# Row 38
$ws.Range("H38").Value = 9989.5
$ws.Range("J38").Value = 9989.5
$ws.Range("L38").Value = 9989.5
$ws.Range("N38").Value = -10915.5
# Row 126
$ws.Range("H126").Value = 2999.5
$ws.Range("I126").Value = 2999.5
$ws.Range("K126").Value = 8998.5
$ws.Range("M126").Value = -6528.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1749
$ws.Range("I40").Value = 1749
$ws.Range("K40").Value = 1749
$ws.Range("M40").Value = -1613
# Row 82
$ws.Range("H82").Value = 4345
$ws.Range("I82").Value = 4650
$ws.Range("K82").Value = 4650
$ws.Range("M82").Value = -4289
# Row 85
$ws.Range("H85").Value = 4345
$ws.Range("I85").Value = 4650
$ws.Range("K85").Value = 4650
$ws.Range("M85").Value = -3402
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
# Row 122
$ws.Range("H122").Value = 4763.8823
$ws.Range("I122").Value = 4135.364
$ws.Range("K122").Value = 12406.092
$ws.Range("M122").Value = -9956.091999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 100000000
$ws.Range("J46").Value = 100000000
$ws.Range("L46").Value = 100000000
$ws.Range("N46").Value = -100000462
# Row 62
$ws.Range("H62").Value = 14000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 14000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 14000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -15248
# Row 65
$ws.Range("H65").Value = 14000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 14000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 70000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -76240
# Row 81
$ws.Range("H81").Value = 4263.593
$ws.Range("I81").Value = 2487.9092
$ws.Range("K81").Value = 4975.8184
$ws.Range("M81").Value = -3914.8184
# Row 84
$ws.Range("H84").Value = 4263.593
$ws.Range("I84").Value = 2487.9092
$ws.Range("K84").Value = 24879.092
$ws.Range("M84").Value = -19575.092
# Row 122
$ws.Range("H122").Value = 5045.375
$ws.Range("I122").Value = 5002
$ws.Range("K122").Value = 15006
$ws.Range("M122").Value = -12556
# Row 126
$ws.Range("H126").Value = 4356.7144
$ws.Range("I126").Value = 3916.5
$ws.Range("K126").Value = 11749.5
$ws.Range("M126").Value = -9279.5
# Row 132
$ws.Range("H132").Value = 6693.278
$ws.Range("I132").Value = 4289.9165
$ws.Range("K132").Value = 12869.7495
$ws.Range("M132").Value = -10339.7495
# Row 134
$ws.Range("H134").Value = 100000000
$ws.Range("J134").Value = 100000000
$ws.Range("L134").Value = 300000000
$ws.Range("N134").Value = -300005070
